$wb = $excel.ActiveWorkbook

# Rename "Red Line" -> "Red Line!" and "Green Line" -> "Green Line!"
$wb.Worksheets.Item("Red Line").Name = "Red Line!"
$wb.Worksheets.Item("Green Line").Name = "Green Line!"

# Shrink the shared formula range for A71 on "Red Line!" from A71:A134 to A71:A77
$ws = $wb.Worksheets.Item("Red Line!")
$ws.Range("A71:A77").Formula = "=A70"
